$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 42, shifting existing rows 42:63 down to 43:64
# (mirrors the source data gaining a new weekly Papaya price record).
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record's data.
$ws.Cells.Item(42, 1).Value = 3
$ws.Cells.Item(42, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(42, 3).Value = "Coquimbo"
$ws.Cells.Item(42, 4).Value = "2023-03-21"
$ws.Cells.Item(42, 5).Value = 5
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100108
$ws.Cells.Item(42, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(42, 9).Value = 100108004
$ws.Cells.Item(42, 10).Value = "Papaya"
$ws.Cells.Item(42, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 48
$ws.Cells.Item(42, 14).Value = 20000
$ws.Cells.Item(42, 15).Value = 20000
$ws.Cells.Item(42, 16).Value = 20000
$ws.Cells.Item(42, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(42, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(42, 19).Value = 2000
$ws.Cells.Item(42, 20).Value = 10
